$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coverage Summary")

$ws.Range("N5").Value = -0.6200000000000001
$ws.Range("O5").Value = -0.6200000000000001
$ws.Range("N6").Value = -0.6200000000000001
$ws.Range("N7").Value = -0.6200000000000001
$ws.Range("N8").Value = -0.6200000000000001
$ws.Range("N21").Value = 0.0
$ws.Range("O21").Value = 0.0
$ws.Range("N25").Value = -0.31000000000000005
$ws.Range("O25").Value = -0.020000000000000462
$ws.Range("N26").Value = -0.2599999999999998
$ws.Range("O26").Value = -0.31000000000000005
$ws.Range("N27").Value = 0.0
$ws.Range("N32").Value = -0.17999999999999994
$ws.Range("O32").Value = -0.7900000000000001
$ws.Range("N33").Value = -2.51
$ws.Range("O33").Value = -0.17999999999999994
$ws.Range("N34").Value = -0.33000000000000007
$ws.Range("O34").Value = -2.37
$ws.Range("N35").Value = 0.0
$ws.Range("O35").Value = 0.0
$ws.Range("N50").Value = -0.21999999999999975
$ws.Range("O50").Value = -0.13000000000000078
$ws.Range("N55").Value = -1.23
$ws.Range("O55").Value = -0.75
$ws.Range("N56").Value = -1.0
$ws.Range("O56").Value = -1.23
$ws.Range("N57").Value = -2.94
$ws.Range("O57").Value = -0.5899999999999999
$ws.Range("N58").Value = 0.0
$ws.Range("O58").Value = -2.69
$ws.Range("O59").Value = 0.0
$ws.Range("N67").Value = -0.040000000000000036
$ws.Range("O67").Value = -0.08999999999999986
$ws.Range("N79").Value = -1.0
$ws.Range("O79").Value = -0.33000000000000007
$ws.Range("N80").Value = -0.8700000000000003
$ws.Range("O80").Value = -0.9999999999999999
$ws.Range("N81").Value = -2.43
$ws.Range("N82").Value = 0.0
$ws.Range("O82").Value = -2.43
$ws.Range("N83").Value = -0.11000000000000032
$ws.Range("O83").Value = 0.0
$ws.Range("N84").Value = -0.03000000000000025
$ws.Range("O84").Value = -0.11000000000000032
$ws.Range("N85").Value = -0.36000000000000026
$ws.Range("O85").Value = -0.03000000000000025
$ws.Range("N97").Value = -0.3200000000000003
$ws.Range("O97").Value = -0.08000000000000007
$ws.Range("N98").Value = -0.08000000000000007
$ws.Range("O98").Value = -0.3200000000000003
$ws.Range("N104").Value = -0.86
$ws.Range("O104").Value = -0.9699999999999999
$ws.Range("N105").Value = -2.48
$ws.Range("O105").Value = 0.0
$ws.Range("N106").Value = 0.0
$ws.Range("O106").Value = -2.48
$ws.Range("N107").Value = 0.0
$ws.Range("O107").Value = 0.0
$ws.Range("N108").Value = 0.0
$ws.Range("N109").Value = -0.3200000000000003
$ws.Range("O109").Value = 0.0
$ws.Range("N121").Value = -0.3100000000000005
$ws.Range("O121").Value = -0.08000000000000007
$ws.Range("N128").Value = -0.51
$ws.Range("O128").Value = -0.30999999999999994
$ws.Range("N129").Value = -1.4300000000000004
$ws.Range("N130").Value = 0.0
$ws.Range("O130").Value = -1.4300000000000004
$ws.Range("O131").Value = 0.0
$ws.Range("O143").Value = 0.0
$ws.Range("N144").Value = -0.1599999999999997
$ws.Range("O144").Value = -0.22000000000000022
$ws.Range("N145").Value = -0.09999999999999964
$ws.Range("N151").Value = -0.2
$ws.Range("O151").Value = -0.28
$ws.Range("N152").Value = -0.10000000000000003
$ws.Range("O152").Value = -0.2
$ws.Range("N153").Value = 0.0
$ws.Range("B174").Value = 29.879999999999992
$ws.Range("B175").Value = 22.4
$ws.Range("N20").Value = 0.0
$ws.Range("O20").Value = -0.34999999999999964
$ws.Range("N142").Value = 0.0
$ws.Range("O142").Value = -0.1499999999999999
